$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 91.7
$ws.Range("D5").Value = 10.494
$ws.Range("D6").Value = 49
$ws.Range("D15").Value = 88.12
$ws.Range("D16").Value = 10.821999999999999
$ws.Range("D17").Value = 55
$ws.Range("D26").Value = 86.32
$ws.Range("D27").Value = 11.175000000000001
$ws.Range("D28").Value = 55
$ws.Range("D37").Value = 85.778999999999996
$ws.Range("D38").Value = 11.5365
$ws.Range("D39").Value = 48
$ws.Range("D44").Value = 2
$ws.Range("D48").Value = 82.88
$ws.Range("D49").Value = 11.929
$ws.Range("D50").Value = 59
$ws.Range("D59").Value = 81.260000000000005
$ws.Range("D60").Value = 12.362500000000001
$ws.Range("D61").Value = 63
$ws.Range("D70").Value = 80.16
$ws.Range("D71").Value = 12.811
$ws.Range("D72").Value = 61
$ws.Range("D81").Value = 72.52
$ws.Range("D82").Value = 10.4895
$ws.Range("D83").Value = 43
$ws.Range("D88").Value = 3
$ws.Range("D92").Value = 71
$ws.Range("D93").Value = 10.81
$ws.Range("D94").Value = 46
$ws.Range("D99").Value = 1.5
$ws.Range("D103").Value = 69.891999999999996
$ws.Range("D104").Value = 11.146000000000001
$ws.Range("D105").Value = 44
$ws.Range("D110").Value = 1.5
$ws.Range("D114").Value = 68.492000000000004
$ws.Range("D115").Value = 11.5045
$ws.Range("D116").Value = 48
$ws.Range("D121").Value = 1.5
$ws.Range("D125").Value = 67.552000000000007
$ws.Range("D126").Value = 11.897500000000001
$ws.Range("D127").Value = 54
$ws.Range("D132").Value = 1.5
$ws.Range("D136").Value = 66.251999999999995
$ws.Range("D137").Value = 12.324999999999999
$ws.Range("D138").Value = 54
$ws.Range("D143").Value = 1.5
$ws.Range("D147").Value = 65.531999999999996
$ws.Range("D148").Value = 12.7735
$ws.Range("D149").Value = 54
$ws.Range("D154").Value = 1.5
$ws.Range("D158").Value = 59.78
$ws.Range("D159").Value = 10.436
$ws.Range("D160").Value = 46
$ws.Range("D165").Value = 2
$ws.Range("D169").Value = 59.472000000000001
$ws.Range("D170").Value = 10.765000000000001
$ws.Range("D171").Value = 40
$ws.Range("D176").Value = 2
$ws.Range("D180").Value = 58.712000000000003
$ws.Range("D181").Value = 11.095000000000001
$ws.Range("D182").Value = 45
$ws.Range("D187").Value = 2
$ws.Range("D191").Value = 57.911999999999999
$ws.Range("D192").Value = 11.4475
$ws.Range("D193").Value = 44
$ws.Range("D198").Value = 2
$ws.Range("D202").Value = 56.792000000000002
$ws.Range("D203").Value = 11.842000000000001
$ws.Range("D204").Value = 48
$ws.Range("D209").Value = 2
$ws.Range("D213").Value = 56.052
$ws.Range("D214").Value = 12.265499999999999
$ws.Range("D215").Value = 49
$ws.Range("D220").Value = 2
$ws.Range("D224").Value = 55.512
$ws.Range("D225").Value = 12.7235
$ws.Range("D226").Value = 47
$ws.Range("D231").Value = 4
$ws.Range("D235").Value = 51.44
$ws.Range("D236").Value = 10.407500000000001
$ws.Range("D237").Value = 40
$ws.Range("D242").Value = 3
$ws.Range("D246").Value = 50.792000000000002
$ws.Range("D247").Value = 10.734999999999999
$ws.Range("D248").Value = 43
$ws.Range("D253").Value = 3
$ws.Range("D257").Value = 50.112000000000002
$ws.Range("D258").Value = 11.081
$ws.Range("D259").Value = 42
$ws.Range("D264").Value = 3
$ws.Range("D268").Value = 49.411999999999999
$ws.Range("D269").Value = 11.423999999999999
$ws.Range("D270").Value = 95
$ws.Range("D275").Value = 3
$ws.Range("D279").Value = 48.731999999999999
$ws.Range("D281").Value = 48
$ws.Range("D286").Value = 2
$ws.Range("D290").Value = 48.152000000000001
$ws.Range("D292").Value = 50
$ws.Range("D297").Value = 2
$ws.Range("D301").Value = 47.832000000000001
$ws.Range("D303").Value = 45
$ws.Range("D308").Value = 2
$ws.Range("D312").Value = 44.795000000000002
$ws.Range("D313").Value = 10.384499999999999
$ws.Range("D314").Value = 40
$ws.Range("D323").Value = 44.351999999999997
$ws.Range("D324").Value = 10.7095
$ws.Range("D325").Value = 38
$ws.Range("D334").Value = 43.832000000000001
$ws.Range("D335").Value = 11.042
$ws.Range("D336").Value = 43

$ws.Range("F338").Select()

